# "Revert "Revert "Merge branch 'master' ... Machine-Learning-Nose-Jobs"""
#
# Net effect on Sheet1:
#   - The backlog row "Apply changes on image of patient face for the use of
#     the surgeon" (old row 10) is removed entirely, shifting every row below
#     it up by one.
#   - The item text "Get ID 7-9 on the python app" is corrected to
#     "Get ID 7-8 on the python app".
#   - Column B is widened and the selection / scroll position are updated to
#     reflect where the author was working afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the row for "Apply changes on image of patient face for the use
#    of the surgeon" (row 10). Rows 11-15 shift up to become rows 10-14.
$ws.Rows.Item(10).Delete() | Out-Null

# 2. The deletion invalidates the formula that used to live in old row 11
#    (now row 10); re-point it so the A-column numbering (1,2,3,...) keeps
#    counting correctly down to row 13.
$ws.Cells.Item(10, 1).Formula = "=A9+1"

# 3. Correct the "Get ID 7-9" story to "Get ID 7-8" (now on row 12).
$ws.Cells.Item(12, 2).Value = "Get ID 7-8 on the python app"

# 4. Column B got widened.
$ws.Columns.Item(2).ColumnWidth = 86.6

# 5. Selection / scroll state left behind by the author.
$ws.Range("A9:A13").Select() | Out-Null

Write-Output "Applied backlog row removal + 'Get ID 7-8' correction"
